# OTreeAlg: confirm (SVS) instead of identify (ADS) for next states of extra states
#
# The workbook's data sheet (Mealy_R100) feeds three charts that live on
# separate chart sheets (TeacherDFSM / TeacherRL / TeacherBB). Those charts
# cache the "ExtraStates:1"/OTree data point (idx 12) from rows 14, 27 and
# 40 of this sheet. The underlying run changed from an ADS-style "identify"
# strategy to an SVS-style "confirm" strategy for extra-state next-state
# resolution, which altered the #Resets / #OQs / #EQs counts recorded for
# that data point — update the three source rows accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (OTree, chart "TeacherDFSM" source range $B$2:$B$14 etc.)
$ws.Range("B14").Value = 2875
$ws.Range("C14").Value = 6040
$ws.Range("E14").Value = 22008

# Row 27 (OTree, chart "TeacherRL" source range $B$15:$B$27 etc.)
$ws.Range("B27").Value = 2875
$ws.Range("C27").Value = 6040
$ws.Range("E27").Value = 22008

# Row 40 (OTree, chart "TeacherBB" source range $B$28:$B$40 etc.)
$ws.Range("B40").Value = 2875
$ws.Range("C40").Value = 6040
$ws.Range("E40").Value = 6040

# The author's last interaction moved the selection off the bottom of the
# sheet (past the last data row) instead of leaving it parked on E40 with
# the view scrolled to row 28.
$ws.Range("A41").Select()
